$d = $word.ActiveDocument

$replacements = @(
    @{old = "485×8=3880"; new = "359×3=1077"},
    @{old = "345×7=2415"; new = "239×7=1673"},
    @{old = "975×7=6825"; new = "647×2=1294"},
    @{old = "550×4=2200"; new = "716×6=4296"},
    @{old = "180×9=1620"; new = "299×3=897"},
    @{old = "754×5=3770"; new = "338×3=1014"},
    @{old = "307×6=1842"; new = "808×2=1616"},
    @{old = "666×8=5328"; new = "731×8=5848"},
    @{old = "891×5=4455"; new = "757×3=2271"},
    @{old = "754×8=6032"; new = "962×6=5772"},
    @{old = "884×3=2652"; new = "813×8=6504"},
    @{old = "965×8=7720"; new = "401×4=1604"},
    @{old = "587×4=2348"; new = "822×4=3288"},
    @{old = "819×2=1638"; new = "239×6=1434"},
    @{old = "799×3=2397"; new = "720×8=5760"},
    @{old = "171×6=1026"; new = "833×9=7497"},
    @{old = "342×8=2736"; new = "762×2=1524"},
    @{old = "842×3=2526"; new = "450×7=3150"},
    @{old = "933×6=5598"; new = "473×7=3311"},
    @{old = "223×2=446"; new = "852×4=3408"},
    @{old = "967×7=6769"; new = "480×7=3360"},
    @{old = "746×5=3730"; new = "766×3=2298"},
    @{old = "920×5=4600"; new = "733×7=5131"},
    @{old = "703×6=4218"; new = "586×2=1172"},
    @{old = "277×9=2493"; new = "885×6=5310"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
